$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (keep stored as text, matching the original t="str" typing).
# Briefly mark the cell as Text so Excel doesn't auto-coerce the numeric-looking
# string into a Number, then restore the default "Normal" style so no extra
# cell formatting is introduced.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "10"
$ws.Range("E2").Style = "Normal"

$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2"
$ws.Range("Y2").Style = "Normal"

$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "0"
$ws.Range("AA2").Style = "Normal"

# Remove row 3 entirely (second data record dropped)
$ws.Rows("3").Delete()
